$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$val111 = @'
Slower than the Medium Tank, but possessing heavier armor and fire power, the Heavy Tank is a fearsome combatant. Trading speed or weapons for armor, the Heavy Tank can soak up a surprisingly heavy level of fire power.
  <b><color=#ffcc00>Wheeled: Fast in open ground, slow in rough terrain.</color></b>
  <b><color=#ffcc00>ICE: Air breathing Combustion Engine. This unit can not deploy in thin atmosphere.</color></b>
  <b><color=#ffcc00>Movement: 3/5 Hex: 90/150 Meters.</color></b>
 <b>Armor:</b> Ferro-Fibrous
 <b>Structure:</b> Structure
 <b>Values A S </b>
 <b>Front</b> 350 35
 <b>Left</b> 282 35
 <b>Right</b> 282 35
 <b>Rear</b> 271 35
 <b>Turret</b> 240 35
 <b>Total</b> 1425 175
---
Slower than the Medium Tank, but possessing heavier armor and fire power, the Heavy Tank is a fearsome combatant. Trading speed or weapons for armor, the Heavy Tank can soak up a surprisingly heavy level of fire power.

 <b><color=#ffcc00>Wheeled: Fast in open ground, slow in rough terrain.</color></b>

 <b><color=#ffcc00>ICE: Air breathing Combustion Engine. This unit can not deploy in thin atmosphere.</color></b>

 <b><color=#ffcc00>Movement: 3/5 Hex: 90/150 Meters.</color></b>

<b>Armor:</b> Ferro-Fibrous
<b>Structure:</b> Structure
<b>Values       A        S  </b>
<b>Front</b>       350     35
<b>Left</b>         282     35
<b>Right</b>       282     35
<b>Rear</b>        247     35
<b>Turret</b>      240     35

<b>Total</b>      1401    175

'@

$val115 = @'
Slower than the Medium Tank, but possessing heavier armor and fire power, the Heavy Tank is a fearsome combatant. Trading speed or weapons for armor, the Heavy Tank can soak up a surprisingly heavy level of fire power.
  <b><color=#ffcc00>Wheeled: Fast in open ground, slow in rough terrain.</color></b>
  <b><color=#ffcc00>Fusion: Expensive engine that works in a Vacuum.</color></b>
  <b><color=#ffcc00>Movement: 3/5 Hex: 90/150 Meters.</color></b>
 <b>Armor:</b> Ferro-Fibrous
 <b>Structure:</b> Structure
 <b>Values A S </b>
 <b>Front</b> 374 35
 <b>Left</b> 269 35
 <b>Right</b> 269 35
 <b>Rear</b> 263 35
 <b>Turret</b> 250 35
 <b>Total</b> 1425 175
---
Slower than the Medium Tank, but possessing heavier armor and fire power, the Heavy Tank is a fearsome combatant. Trading speed or weapons for armor, the Heavy Tank can soak up a surprisingly heavy level of fire power.

 <b><color=#ffcc00>Wheeled: Fast in open ground, slow in rough terrain.</color></b>

 <b><color=#ffcc00>Fusion: Expensive engine that works in a Vacuum.</color></b>

 <b><color=#ffcc00>Movement: 3/5 Hex: 90/150 Meters.</color></b>

<b>Armor:</b> Ferro-Fibrous
<b>Structure:</b> Structure
<b>Values       A        S  </b>
<b>Front</b>       374     35
<b>Left</b>         269     35
<b>Right</b>       269     35
<b>Rear</b>        239     35
<b>Turret</b>      250     35

<b>Total</b>      1401    175

'@

$ws.Range("B51").Value = $val111
$ws.Range("B76").Value = $val111
$ws.Range("B53").Value = $val115
$ws.Range("B77").Value = $val115
